$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2998339386254827
$ws.Range("C2").Value = 0.04903159714902472
$ws.Range("D2").Value = 0.07851510717721055
$ws.Range("E2").Value = 0.1494589360147245
$ws.Range("G2").Value = 0.002444180571458521
$ws.Range("I2").Value = 0.6398288474074398
$ws.Range("K2").Value = 0.3237104454460109
$ws.Range("M2").Value = 0.2401359665607004
$ws.Range("N2").Value = 1.593526331344803
$ws.Range("O2").Value = 3.040280488031499

$ws.Range("B3").Value = 0.2673024519064882
$ws.Range("C3").Value = 0.04277111620473306
$ws.Range("D3").Value = 0.07122575567662182
$ws.Range("E3").Value = 0.1386873054833515
$ws.Range("G3").Value = 0.002446796989621838
$ws.Range("I3").Value = 0.6427850256888981
$ws.Range("K3").Value = 0.2865146929551656
$ws.Range("M3").Value = 0.2175135535620853
$ws.Range("N3").Value = 1.610609169924331
$ws.Range("O3").Value = 3.041201842892008

$ws.Range("B4").Value = 0.2473687087033341
$ws.Range("C4").Value = 0.03892031210686753
$ws.Range("D4").Value = 0.06678392970682978
$ws.Range("E4").Value = 0.1321716490359108
$ws.Range("G4").Value = 0.002448489620338353
$ws.Range("I4").Value = 0.6449280029245053
$ws.Range("K4").Value = 0.2637016229586777
$ws.Range("M4").Value = 0.2037189408445528
$ws.Range("N4").Value = 1.62164069040888
$ws.Range("O4").Value = 3.043449654035072

$ws.Range("B5").Value = 0.2392561369210284
$ws.Range("C5").Value = 0.03734937123758186
$ws.Range("D5").Value = 0.06498237364364456
$ws.Range("E5").Value = 0.1295409430060062
$ws.Range("G5").Value = 0.002449201106830088
$ws.Range("I5").Value = 0.6458837030429692
$ws.Range("K5").Value = 0.2544117938443691
$ws.Range("M5").Value = 0.1981214969055216
$ws.Range("N5").Value = 1.62627250032974
$ws.Range("O5").Value = 3.04478846843719

$ws.Range("B6").Value = 0.2379097005405413
$ws.Range("C6").Value = 0.03708841502886173
$ws.Range("D6").Value = 0.06468374171632263
$ws.Range("E6").Value = 0.1291055885368166
$ws.Range("G6").Value = 0.002449320562701333
$ws.Range("I6").Value = 0.6460473734912284
$ws.Range("K6").Value = 0.252869638044416
$ws.Range("M6").Value = 0.1971934912876989
$ws.Range("N6").Value = 1.627049843807617
$ws.Range("O6").Value = 3.04503630972269

$ws.Range("B7").Value = 0.2472592563779585
$ws.Range("C7").Value = 0.03889913274134926
$ws.Range("D7").Value = 0.06675959876621107
$ws.Range("E7").Value = 0.1321360715831901
$ws.Range("G7").Value = 0.002448499127623364
$ws.Range("I7").Value = 0.6449405581534506
$ws.Range("K7").Value = 0.2635763095840815
$ws.Range("M7").Value = 0.2036433547937833
$ws.Range("N7").Value = 1.621702604522136
$ws.Range("O7").Value = 3.043465998022413

$ws.Range("B8").Value = 0.2886088305506576
$ws.Range("C8").Value = 0.04687441016001515
$ws.Range("D8").Value = 0.07599471281724846
$ws.Range("E8").Value = 0.1457243741718628
$ws.Range("G8").Value = 0.002445064873695708
$ws.Range("I8").Value = 0.6407800710096154
$ws.Range("K8").Value = 0.3108802894811333
$ws.Range("M8").Value = 0.2323158550491371
$ws.Range("N8").Value = 1.599303843773807
$ws.Range("O8").Value = 3.040248945088337

$ws.Range("B9").Value = 0.370006985792287
$ws.Range("C9").Value = 0.06245975833662953
$ws.Range("D9").Value = 0.09437405003346555
$ws.Range("E9").Value = 0.1731607473884154
$ws.Range("G9").Value = 0.00243901080189496
$ws.Range("I9").Value = 0.6352245235481035
$ws.Range("K9").Value = 0.4038340221846681
$ws.Range("M9").Value = 0.2893078667147293
$ws.Range("N9").Value = 1.559687262972691
$ws.Range("O9").Value = 3.047298375917478

$ws.Range("B10").Value = 0.4299910426904603
$ws.Range("C10").Value = 0.07387864246659603
$ws.Range("D10").Value = 0.1080437660276772
$ws.Range("E10").Value = 0.1938167660784202
$ws.Range("G10").Value = 0.002434973534653343
$ws.Range("I10").Value = 0.6327327257882018
$ws.Range("K10").Value = 0.4722371454320182
$ws.Range("M10").Value = 0.3316597282936513
$ws.Range("N10").Value = 1.533208760437635
$ws.Range("O10").Value = 3.060642588897963

$ws.Range("B11").Value = 0.4573170159726203
$ws.Range("C11").Value = 0.0790668653651494
$ws.Range("D11").Value = 0.1142991657434607
$ws.Range("E11").Value = 0.2033257639171993
$ws.Range("G11").Value = 0.002433225157634708
$ws.Range("I11").Value = 0.631945047430122
$ws.Range("K11").Value = 0.503378466914711
$ws.Range("M11").Value = 0.3510339802253313
$ws.Range("N11").Value = 1.52173403222546
$ws.Range("O11").Value = 3.068491230031555

$ws.Range("B12").Value = 0.4676699749320221
$ws.Range("C12").Value = 0.08103060943670926
$ws.Range("D12").Value = 0.1166732516331734
$ws.Range("E12").Value = 0.2069430074005396
$ws.Range("G12").Value = 0.002432575708134769
$ws.Range("I12").Value = 0.6316965547110343
$ws.Range("K12").Value = 0.5151741463830604
$ws.Range("M12").Value = 0.3583862130640583
$ws.Range("N12").Value = 1.517470960900926
$ws.Range("O12").Value = 3.071719332744294

$ws.Range("B13").Value = 0.4654400546790782
$ws.Range("C13").Value = 0.08060772339408118
$ws.Range("D13").Value = 0.1161617142152238
$ws.Range("E13").Value = 0.2061632367485871
$ws.Range("G13").Value = 0.002432715018143964
$ws.Range("I13").Value = 0.6317478570965136
$ws.Range("K13").Value = 0.5126336021602356
$ws.Range("M13").Value = 0.3568020822039202
$ws.Range("N13").Value = 1.51838543350328
$ws.Range("O13").Value = 3.071012713308107

$ws.Range("B14").Value = 0.4581686580444853
$ws.Range("C14").Value = 0.07922844238972004
$ws.Range("D14").Value = 0.1144943770542568
$ws.Range("E14").Value = 0.2036230269318366
$ws.Range("G14").Value = 0.002433171474380819
$ws.Range("I14").Value = 0.6319236057477013
$ws.Range("K14").Value = 0.5043488438241752
$ws.Range("M14").Value = 0.3516385393254922
$ws.Range("N14").Value = 1.521381659459646
$ws.Range("O14").Value = 3.06875167592321

$ws.Range("B15").Value = 0.4537153839971495
$ws.Range("C15").Value = 0.07838347212052099
$ws.Range("D15").Value = 0.1134737748757715
$ws.Range("E15").Value = 0.2020692163510418
$ws.Range("G15").Value = 0.002433452709275279
$ws.Range("I15").Value = 0.6320377417263501
$ws.Range("K15").Value = 0.4992745876128311
$ws.Range("M15").Value = 0.348477756959177
$ws.Range("N15").Value = 1.523227637701616
$ws.Range("O15").Value = 3.067400071262739

$ws.Range("B16").Value = 0.428205975880843
$ws.Range("C16").Value = 0.07353945175904641
$ws.Range("D16").Value = 0.1076357037113382
$ws.Range("E16").Value = 0.1931976126442052
$ws.Range("G16").Value = 0.002435089564893894
$ws.Range("I16").Value = 0.6327911661974994
$ws.Range("K16").Value = 0.4702024460265761
$ws.Range("M16").Value = 0.330395758023478
$ws.Range("N16").Value = 1.533970151002503
$ws.Range("O16").Value = 3.06016546912349

$ws.Range("B17").Value = 0.4125664713536992
$ws.Range("C17").Value = 0.07056617854149749
$ws.Range("D17").Value = 0.1040636960907335
$ws.Range("E17").Value = 0.1877841335710499
$ws.Range("G17").Value = 0.00243611626926889
$ws.Range("I17").Value = 0.6333419849524802
$ws.Range("K17").Value = 0.4523736133632212
$ws.Range("M17").Value = 0.3193307929564071
$ws.Range("N17").Value = 1.54070653582712
$ws.Range("O17").Value = 3.056182934568767

$ws.Range("B18").Value = 0.4035747291972029
$ws.Range("C18").Value = 0.06885544004776989
$ws.Range("D18").Value = 0.102012653235505
$ws.Range("E18").Value = 0.1846810168153539
$ws.Range("G18").Value = 0.002436715107451558
$ws.Range("I18").Value = 0.6336913515880696
$ws.Range("K18").Value = 0.4421212590054324
$ws.Range("M18").Value = 0.3129766958493292
$ws.Range("N18").Value = 1.544634814350591
$ws.Range("O18").Value = 3.054059650631558

$ws.Range("B19").Value = 0.40053092500969
$ws.Range("C19").Value = 0.06827611258727018
$ws.Range("D19").Value = 0.1013188037712354
$ws.Range("E19").Value = 0.183632164088074
$ws.Range("G19").Value = 0.002436919291906695
$ws.Range("I19").Value = 0.6338152301268991
$ws.Range("K19").Value = 0.4386503975239009
$ws.Range("M19").Value = 0.3108270537646902
$ws.Range("N19").Value = 1.545974079178524
$ws.Range("O19").Value = 3.053369479812091

$ws.Range("B20").Value = 0.414230944534097
$ws.Range("C20").Value = 0.07088274979508924
$ws.Range("D20").Value = 0.1044435821936816
$ws.Range("E20").Value = 0.1883593117857743
$ws.Range("G20").Value = 0.002436006115709784
$ws.Range("I20").Value = 0.6332799801581501
$ws.Range("K20").Value = 0.4542712857575282
$ws.Range("M20").Value = 0.3205076234836497
$ws.Range("N20").Value = 1.539983878736734
$ws.Range("O20").Value = 3.056589558555572

$ws.Range("B21").Value = 0.4603043053298279
$ws.Range("C21").Value = 0.07963359587080276
$ws.Range("D21").Value = 0.1149839706335882
$ws.Range("E21").Value = 0.204368701498737
$ws.Range("G21").Value = 0.002433037060031075
$ws.Range("I21").Value = 0.6318706325796555
$ws.Range("K21").Value = 0.5067821960002163
$ws.Range("M21").Value = 0.3531547724587867
$ws.Range("N21").Value = 1.520499364474261
$ws.Range("O21").Value = 3.069408847635088

$ws.Range("B22").Value = 0.490446072498969
$ws.Range("C22").Value = 0.0853473886122913
$ws.Range("D22").Value = 0.1219036189608289
$ws.Range("E22").Value = 0.2149274123265599
$ws.Range("G22").Value = 0.002431170158051181
$ws.Range("I22").Value = 0.6312397453789913
$ws.Range("K22").Value = 0.5411192230340305
$ws.Range("M22").Value = 0.3745826722901882
$ws.Range("N22").Value = 1.508244055105752
$ws.Range("O22").Value = 3.079279244334856

$ws.Range("B23").Value = 0.474356213558849
$ws.Range("C23").Value = 0.08229832964914863
$ws.Range("D23").Value = 0.1182076508030576
$ws.Range("E23").Value = 0.2092832052638443
$ws.Range("G23").Value = 0.002432159848695351
$ws.Range("I23").Value = 0.6315498913172206
$ws.Range("K23").Value = 0.5227913710408814
$ws.Range("M23").Value = 0.3631378406799257
$ws.Range("N23").Value = 1.514741078458634
$ws.Range("O23").Value = 3.073874595219053

$ws.Range("B24").Value = 0.4134784373724756
$ws.Range("C24").Value = 0.07073963217963808
$ws.Range("D24").Value = 0.1042718277370795
$ws.Range("E24").Value = 0.188099245188134
$ws.Range("G24").Value = 0.002436055889313941
$ws.Range("I24").Value = 0.6333079106888064
$ws.Range("K24").Value = 0.453413355271266
$ws.Range("M24").Value = 0.3199755556764643
$ws.Range("N24").Value = 1.540310419313808
$ws.Range("O24").Value = 3.056405205759347

$ws.Range("B25").Value = 0.3479542059584162
$ws.Range("C25").Value = 0.05824916721438456
$ws.Range("D25").Value = 0.08937288747894456
$ws.Range("E25").Value = 0.1656521703624634
$ws.Range("G25").Value = 0.002440576171938072
$ws.Range("I25").Value = 0.6364484614780466
$ws.Range("K25").Value = 0.3786678008593185
$ws.Range("M25").Value = 0.2738067115909004
$ws.Range("N25").Value = 1.569943280120494
$ws.Range("O25").Value = 3.043959023283094
